# Generate Report for handoff
# Updates the "Latest Handoff Datetime" for the row that is being handed off
# (source file 4c4de13b-...) on both the zh-cn and de-de localization sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-25 13:22:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-25 13:22:33"
